# Updates the "cryptos" sheet with refreshed price/volume figures (and a
# few re-ranked coin rows) as captured by the GitHub Actions scraper run.
# Values that look like plain decimal numbers are written with a leading
# apostrophe so Excel keeps them as literal text (matching the original
# text-formatted price column) instead of silently converting them to
# numeric cells and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.512.77'
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").Value = '3.945.54'
$ws.Range("E3").Value = '  +4.50%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''470.91'
$ws.Range("E5").Value = '  +7.54%  '
$ws.Range("D6").Value = '''148.28'
$ws.Range("E6").Value = '  +5.65%  '
$ws.Range("D7").Value = '''0.627'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''0.733'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +9.71%  '
$ws.Range("D11").Value = '''0.0000353'
$ws.Range("E11").Value = '  +10.97%  '
$ws.Range("D12").Value = '''43.40'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '4.574.24'
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").Value = '''10.41'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '''15.08'
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '3.939.78'
$ws.Range("E16").Value = '  +4.91%  '
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '''19.98'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").Value = '67.697.75'
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("E21").Value = '  +5.11%  '
$ws.Range("D22").Value = '''3.40'
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("D23").Value = '''14.50'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '''87.70'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("E25").Value = '  +6.52%  '
$ws.Range("D26").Value = '''38.64'
$ws.Range("E26").Value = '  +4.23%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").Value = '''10.21'
$ws.Range("E27").Value = '  +4.76%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''9.88'
$ws.Range("E28").Value = '  +4.51%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '''723.50'
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.133'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '''13.44'
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '''2.81'
$ws.Range("E32").Value = '  +4.18%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '''42.37'
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").Value = '0.0₃0841'
$ws.Range("E34").Value = '  +24.10%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '''58.07'
$ws.Range("E35").Value = '  +2.78%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.151'
$ws.Range("E36").Value = '  -2.36%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''5.37'
$ws.Range("E38").Value = '  -3.79%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0476'
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").Value = '''3.05'
$ws.Range("E40").Value = '  +5.03%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '''0.142'
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = '''0.336'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '''2.83'
$ws.Range("E44").Value = '  +6.50%  '
$ws.Range("D45").Value = '''3.49'
$ws.Range("E45").Value = '  +5.88%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''2.21'
$ws.Range("E46").Value = '  +6.38%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '''2.54'
$ws.Range("E47").Value = '  -3.18%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''3.28'
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''147.15'
$ws.Range("E49").Value = '  +3.01%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '''2.89'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''25.55'
$ws.Range("E51").Value = '  +3.88%  '
